# Inventario de Materiales Equipos Vehiculos Piloto.xlsx
# "PowerAppsInfo" sheet: rename/repurpose the last two data columns.
#   B: Quantity          -> ExpectedQuantity        (values unchanged)
#   C: AccountCategory    -> AccountCategory         (unchanged)
#   D: Country / "Peru"   -> Active / TRUE           (every data row)
#   E: Active / "Partes.."-> ItemCategory / "Herramienta" (or "Consumible"
#                                                          for the two
#                                                          consumable rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PowerAppsInfo")

# --- header row -------------------------------------------------------
$ws.Range("B1").Value = "ExpectedQuantity"
$ws.Range("D1").Value = "Active"
$ws.Range("E1").Value = "ItemCategory"

# D1 used to carry the "right edge of the box" border (it used to be the
# last column of the bordered header); now that E1 is a real column too,
# D1 drops that border and matches E1's plain (borderless) header style.
$ws.Range("D1").Borders.LineStyle = 0

# --- data rows 2-66 -----------------------------------------------------
$lastRow = 66
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "Herramienta"
}

# Consumable items (nails / wall plugs) get "Consumible" instead of
# "Herramienta" in the new ItemCategory column.
$ws.Cells.Item(15, 5).Value = "Consumible"
$ws.Cells.Item(64, 5).Value = "Consumible"

# --- cosmetic column widths (best-fit recompute after the text/content
#     of columns B-E changed) -------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14.330729166666666
$ws.Columns.Item(3).ColumnWidth = 27.330729166666668
$ws.Columns.Item(4).ColumnWidth = 5.166666666666667

# --- selection moved to C11 in the saved view --------------------------
$ws.Range("C11").Select() | Out-Null
